$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.895.06"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.70"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.79"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3885"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.59"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.332"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08401"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.72"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.986"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.823"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001310"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.644.95"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.80"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06958"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.49"
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.889"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.62"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.901.07"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.906"
$ws.Range("E26").Value = "  -8.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.86"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.10"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.521"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.70"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.504"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.613"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.835.95"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08033"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9787"
$ws.Range("E35").Value = "  -7.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02901"
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.617"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2674"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.36"
$ws.Range("E39").Value = "  -8.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09084"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7502"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.30"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.416"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.50"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6902"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.415"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.084"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08253"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.37"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.218"
$ws.Range("E51").Value = "  -2.41%  "
